# Rules.Main.xlsx: update the "Good Morning" string used in E8 to "GIT UPDATE"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"

# Mirror the active-cell selection on E8 that was captured when the file was saved
$ws.Range("E8").Select()
